# Generate Report for Handback
# Updates the "latest" handoff/handback timestamps recorded on the
# Overview / zh-cn / de-de sheets as part of (re)generating the report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G2 - "Latest HO Xliff Generate Date"
$wsOverview.Range("G2").Value = "2016-09-02 01:16:22"

# zh-cn!H2 - "Correspond Handoff Datetime"
$wsZhCn.Range("H2").Value = "2016-09-02 01:16:17"
# zh-cn!K2 - "Correspond Handback DateTime"
$wsZhCn.Range("K2").Value = "2016-09-02 01:16:41"

# de-de!H2 - "Correspond Handoff Datetime"
$wsDeDe.Range("H2").Value = "2016-09-02 01:16:22"
# de-de!K2 - "Correspond Handback DateTime"
$wsDeDe.Range("K2").Value = "2016-09-02 01:16:49"
